$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the new I1:J1 header cells, then set their text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF data columns for rows 2-11
$iValues = @(6, 7, 5, 3, 8, 7, 7, 6, 8, 7)
$jValues = @(6, 7, 5, 5, 9, 7, 7, 6, 8, 7)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
